# ESC-ISR.xlsx: "Working version with power adjusted control"
#
# Rows 210-226 held a block of per-channel latency samples (START/OFFA/
# LOWA/... labels in column D, stop-watch "Actual" readings in E, OCR
# counts in G, with F/H holding the derived deltas). The new working
# version replaces that block with a blank A/B/C template (just the three
# new labels in D210:D212) and clears out the old sample data in
# E/G210:226 so the dependent F/H formulas settle back to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row labels for the reworked table - three fresh shared strings
# (A/B/C) replace the old START/OFFA/LOWA/... cycle in D210:D212, and
# D213:D226 lose their labels entirely (no cell at all, matching the
# cleared-out template).
$ws.Range("D210").Value = "A"
$ws.Range("D211").Value = "B"
$ws.Range("D212").Value = "C"
$ws.Range("D213:D226").ClearContents()

# Old stop-watch "Actual" readings in E210:E211 and OCR @Case counts in
# G210:G226 are cleared; dependent formulas in F211:F226/H210:H226
# recalculate to 0 automatically.
$ws.Range("E210:E211").ClearContents()
$ws.Range("G210:G226").ClearContents()

# Scroll the saved selection up from E212 to E210 to match where the
# cursor was left.
$ws.Range("E210").Select() | Out-Null
